# Update countries & provincias Spain
# Applies updated COVID case numbers and re-sorted country rows
# (as some countries changed rank / total cases, swapping adjacent rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Agosto de 2020 a las 05:54"

# --- Row 22: Alemania (no reordering, only Casos activos / Recuperados change) ---
$ws.Range("D22").Value = 196100
$ws.Range("E22").Value = 8759

# --- Row 29: Kazajistan (no reordering) ---
$ws.Range("B29").Value = 95942
$ws.Range("C29").Value = 1060
$ws.Range("D29").Value = 68871
$ws.Range("E29").Value = 26013

# --- Row 39: Belgica (no reordering) ---
$ws.Range("B39").Value = 71158
$ws.Range("C39").Value = 510
$ws.Range("D39").Value = 17661
$ws.Range("E39").Value = 43638
$ws.Range("G39").Value = 7
$ws.Range("H39").Value = 9859

# --- Rows 50-51: Nigeria / Honduras swap positions (Honduras moves up to 50) ---
$ws.Range("A50").Value = "Honduras"
$ws.Range("B50").Value = 45098
$ws.Range("C50").Value = 799
$ws.Range("D50").Value = 6116
$ws.Range("E50").Value = 37559
$ws.Range("G50").Value = 23
$ws.Range("H50").Value = 1423

$ws.Range("A51").Value = "Nigeria"
$ws.Range("B51").Value = 44890
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 32165
$ws.Range("E51").Value = 11798
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 927

# --- Row 93: Haiti (no reordering) ---
$ws.Range("B93").Value = 7544
$ws.Range("C93").Value = 12
$ws.Range("E93").Value = 2541

# --- Rows 153-154: Malta / Jamaica swap positions (Jamaica moves up to 153) ---
$ws.Range("A153").Value = "Jamaica"
$ws.Range("B153").Value = 928
$ws.Range("C153").Value = 8
$ws.Range("D153").Value = 745
$ws.Range("E153").Value = 171
$ws.Range("H153").Value = 12

$ws.Range("A154").Value = "Malta"
$ws.Range("B154").Value = 926
$ws.Range("C154").Value = 0
$ws.Range("D154").Value = 668
$ws.Range("E154").Value = 249
$ws.Range("H154").Value = 9

# --- Row 160: Vietnam (no reordering) ---
$ws.Range("E160").Value = 327
$ws.Range("G160").Value = 1
$ws.Range("H160").Value = 9

# --- Row 172: Mongolia (no reordering) ---
$ws.Range("D172").Value = 260
$ws.Range("E172").Value = 33

# --- Rows 187-189: Seychelles / Monaco / Islas Turcas y Caicos reordered
#     (Islas Turcas y Caicos moves up to 187) ---
$ws.Range("A187").Value = "Islas Turcas y Caicos"
$ws.Range("B187").Value = 129
$ws.Range("C187").Value = 13
$ws.Range("D187").Value = 39
$ws.Range("E187").Value = 88
$ws.Range("H187").Value = 2

$ws.Range("A188").Value = "Seychelles"
$ws.Range("B188").Value = 126
$ws.Range("C188").Value = 0
$ws.Range("D188").Value = 124
$ws.Range("E188").Value = 2
$ws.Range("H188").Value = 0

$ws.Range("A189").Value = "Monaco"
$ws.Range("B189").Value = 125
$ws.Range("C189").Value = 0
$ws.Range("D189").Value = 105
$ws.Range("E189").Value = 16
$ws.Range("H189").Value = 4

# --- Row 193: Belice (no reordering) ---
$ws.Range("B193").Value = 86
$ws.Range("C193").Value = 14
$ws.Range("E193").Value = 53

# --- Rows 202-203: Timor Oriental / Santa Lucia swap positions
#     (identical figures, so only the country names swap) ---
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Timor Oriental"
